$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Control de excepciones: update delivery status for entry 1244 (row 8) to "Entregado"
$ws.Range("F8").Value = "Entregado"

# Update the last active selection to reflect where work left off
$ws.Range("H9").Select()
